$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 70: hours bumped from 1 to 1.25, and the note changed from
# "2 small problems" to "3 small problems" (combining the two lists).
$ws.Range("C70").Value = 1.25
$ws.Range("D70").Value = "3 small problems"

# Reflect the user's last active cell selection on the sheet.
$ws.Range("C71").Select()
